{"js": "// The sentence \"... avec des maxima de d\u00e9penses fix\u00e9s \u00e0 l'avance pour\n// plusieurs ann\u00e9es. ...\" gets the word \"discr\u00e9tionnaires\" (plus the\n// space that used to separate \"d\u00e9penses\" and \"fix\u00e9s\") inserted right\n// after \"d\u00e9penses\", turning it into\n// \"... avec des maxima de d\u00e9penses discr\u00e9tionnaires fix\u00e9s \u00e0 l'avance ...\".\n//\n// \"maxima de d\u00e9penses \" (with the trailing space) is unique in the\n// document, so searching for it gives us an unambiguous anchor; we then\n// collapse that hit to its end (the caret right before \"fix\u00e9s\") and type\n// the new word plus a single trailing space there. Word/Office.js keeps\n// the insertion point's surrounding run formatting (Times New Roman,\n// 11.5pt/sz 23, fr-FR), so the new text automatically matches the rest\n// of the sentence.\nconst body = context.document.body;\nconst hits = body.search(\"maxima de d\u00e9penses \", { matchCase: true, matchWholeWord: false });\nhits.load(\"text\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the target sentence fragment 'maxima de d\u00e9penses '.\");\n}\n\nconst target = hits.items[0];\nconst insertionPoint = target.getRange(\"End\");\ninsertionPoint.insertText(\"discr\u00e9tionnaires \", \"Before\");\nawait context.sync();\n", "ps1": "# The sentence \"... avec des maxima de d\u00e9penses fix\u00e9s \u00e0 l'avance pour\n# plusieurs ann\u00e9es. ...\" gains the word \"discr\u00e9tionnaires\" right after\n# \"d\u00e9penses\", becoming\n# \"... avec des maxima de d\u00e9penses discr\u00e9tionnaires fix\u00e9s \u00e0 l'avance ...\".\n#\n# \"maxima de d\u00e9penses fix\u00e9s\" is unique in the document, so a plain\n# Find/Replace on that exact (case-sensitive) phrase unambiguously\n# targets the right spot; Word keeps the existing run formatting\n# (Times New Roman, sz 23, fr-FR) for the replacement text since it is\n# typed in place of the matched range.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"maxima de d\u00e9penses fix\u00e9s\"\n$find.Replacement.Text = \"maxima de d\u00e9penses discr\u00e9tionnaires fix\u00e9s\"\n$find.Forward = $true\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Wrap = 0  # wdFindStop: do not wrap around, only the single expected hit\n\n$found = $find.Execute(\n    $find.Text,\n    $false,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    0,\n    $false,\n    $find.Replacement.Text,\n    2  # wdReplaceAll\n)\n\nif (-not $found) {\n    throw \"Could not find the target phrase 'maxima de d\u00e9penses fix\u00e9s' to replace.\"\n}\n"}
